$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 428, pushing existing rows 428-521 down to 430-523
$ws.Rows.Item(428).Resize(2,1).Insert()

# Populate new row 428
$ws.Range("A428").Value2 = 6
$ws.Range("B428").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C428").Value2 = "Metropolitana"
$ws.Range("D428").Value2 = 44694
$ws.Range("E428").Value2 = 13
$ws.Range("F428").Value2 = 100112013
$ws.Range("G428").Value2 = "Alcachofa"
$ws.Range("H428").Value2 = "Española"
$ws.Range("I428").Value2 = "Primera"
$ws.Range("J428").Value2 = 400
$ws.Range("K428").Value2 = 18000
$ws.Range("L428").Value2 = 20000
$ws.Range("M428").Value2 = 18850
$ws.Range("N428").Value2 = "`$/caja 30 unidades"
$ws.Range("O428").Value2 = "Provincia de Limarí"
$ws.Range("P428").Value2 = 628
$ws.Range("Q428").Value2 = 30
$ws.Range("R428").Value2 = "Hortaliza"

# Populate new row 429
$ws.Range("A429").Value2 = 6
$ws.Range("B429").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C429").Value2 = "Metropolitana"
$ws.Range("D429").Value2 = 44694
$ws.Range("E429").Value2 = 13
$ws.Range("F429").Value2 = 100112013
$ws.Range("G429").Value2 = "Alcachofa"
$ws.Range("H429").Value2 = "Española"
$ws.Range("I429").Value2 = "Segunda"
$ws.Range("J429").Value2 = 400
$ws.Range("K429").Value2 = 16000
$ws.Range("L429").Value2 = 18000
$ws.Range("M429").Value2 = 16850
$ws.Range("N429").Value2 = "`$/caja 40 unidades"
$ws.Range("O429").Value2 = "Provincia de Limarí"
$ws.Range("P429").Value2 = 421
$ws.Range("Q429").Value2 = 40
$ws.Range("R429").Value2 = "Hortaliza"
